$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B53: was stored as an inline string "2", should become a real numeric value 2
$ws.Range("B53").Value = 2

# Add new row 54
$ws.Range("A54").Value = "Sunsi Wu"

# B54 keeps its "4" as text (not numeric) like it was scraped into the sheet
$ws.Range("B54").NumberFormat = "@"
$ws.Range("B54").Value = "4"
$ws.Range("B54").Style = "Normal"

$ws.Range("C54").Value = "无"
$ws.Range("D54").Value = "DFT"
$ws.Range("E54").Value = "THE"
$ws.Range("F54").Value = "a0a400ab-cd67-43a0-98e0-d641a379b0a8"
$ws.Range("G54").Value = "B1QRgziT-_annotated.xlsx"
$ws.Range("H54").Value = "I am also interested to hear more about the semantics of the spectral norm of this object (flattened filterbank), which Ian asked about below."
